$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Re-baseline" rows appended after the existing data (rows 27-36).
# Most rows only carry a label in column A, with B:D left blank (but still
# present as real, formatted cells) to match the layout of the rows above.

$ws.Range("A27").Value = "Re-baseline this quarter"
$ws.Range("B27:D27").NumberFormat = "General"

$ws.Range("A28").Value = "Re-baseline ALB/Programme milestones"
$ws.Range("B28:D28").NumberFormat = "General"

$ws.Range("A29").Value = "Re-baseline ALB/Programme cost"
$ws.Range("B29:D29").NumberFormat = "General"

$ws.Range("A30").Value = "Re-baseline ALB/Programme benefits"
$ws.Range("B30:D30").NumberFormat = "General"

$ws.Range("A31").Value = "Re-baseline IPDC milestones"
$ws.Range("B31").Value = "Yes"
$ws.Range("C31").Value = "Yes"
$ws.Range("D31").Value = "Yes"
$ws.Range("E31").Value = "Yes"
$ws.Range("F31").Value = "Yes"

$ws.Range("A32").Value = "Re-baseline IPDC cost"
$ws.Range("B32:D32").NumberFormat = "General"

$ws.Range("A33").Value = "Re-baseline IPDC benefits"
$ws.Range("B33:D33").NumberFormat = "General"

$ws.Range("A34").Value = "Re-baseline HMT milestones"
$ws.Range("B34:D34").NumberFormat = "General"

$ws.Range("A35").Value = "Re-baseline HMT cost"
$ws.Range("B35:D35").NumberFormat = "General"

$ws.Range("A36").Value = "Re-baseline HMT benefits"
$ws.Range("B36:D36").NumberFormat = "General"

# Update the saved selection on the sheet view
[void]$ws.Range("J27").Select()
